$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-03 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-04 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("10+70=80", $true, $false, $false, $false, $false, $true, 1, $false, "80-12=68", 2) | Out-Null
$d.Content.Find.Execute("90-26=64", $true, $false, $false, $false, $false, $true, 1, $false, "71+3=74", 2) | Out-Null
$d.Content.Find.Execute("32+14=46", $true, $false, $false, $false, $false, $true, 1, $false, "55-49=6", 2) | Out-Null
$d.Content.Find.Execute("9+49=58", $true, $false, $false, $false, $false, $true, 1, $false, "5+85=90", 2) | Out-Null
$d.Content.Find.Execute("11+35=46", $true, $false, $false, $false, $false, $true, 1, $false, "86-60=26", 2) | Out-Null
$d.Content.Find.Execute("10+46=56", $true, $false, $false, $false, $false, $true, 1, $false, "71+2=73", 2) | Out-Null
$d.Content.Find.Execute("35-21=14", $true, $false, $false, $false, $false, $true, 1, $false, "20+33=53", 2) | Out-Null
$d.Content.Find.Execute("14+44=58", $true, $false, $false, $false, $false, $true, 1, $false, "78+9=87", 2) | Out-Null
$d.Content.Find.Execute("12+62=74", $true, $false, $false, $false, $false, $true, 1, $false, "87-51=36", 2) | Out-Null
$d.Content.Find.Execute("54-44=10", $true, $false, $false, $false, $false, $true, 1, $false, "22+49=71", 2) | Out-Null
$d.Content.Find.Execute("48+21=69", $true, $false, $false, $false, $false, $true, 1, $false, "56-55=1", 2) | Out-Null
$d.Content.Find.Execute("86-33=53", $true, $false, $false, $false, $false, $true, 1, $false, "94-70=24", 2) | Out-Null
$d.Content.Find.Execute("73-43=30", $true, $false, $false, $false, $false, $true, 1, $false, "63-48=15", 2) | Out-Null
$d.Content.Find.Execute("22-22=0", $true, $false, $false, $false, $false, $true, 1, $false, "91-25=66", 2) | Out-Null
$d.Content.Find.Execute("68-54=14", $true, $false, $false, $false, $false, $true, 1, $false, "77-8=69", 2) | Out-Null
$d.Content.Find.Execute("58-26=32", $true, $false, $false, $false, $false, $true, 1, $false, "48+2=50", 2) | Out-Null
$d.Content.Find.Execute("53-13=40", $true, $false, $false, $false, $false, $true, 1, $false, "18+34=52", 2) | Out-Null
$d.Content.Find.Execute("42-5=37", $true, $false, $false, $false, $false, $true, 1, $false, "78-19=59", 2) | Out-Null
$d.Content.Find.Execute("66-61=5", $true, $false, $false, $false, $false, $true, 1, $false, "40-23=17", 2) | Out-Null
$d.Content.Find.Execute("35+63=98", $true, $false, $false, $false, $false, $true, 1, $false, "31+34=65", 2) | Out-Null
$d.Content.Find.Execute("97-46=51", $true, $false, $false, $false, $false, $true, 1, $false, "44-7=37", 2) | Out-Null
$d.Content.Find.Execute("10+86=96", $true, $false, $false, $false, $false, $true, 1, $false, "86-76=10", 2) | Out-Null
$d.Content.Find.Execute("19+13=32", $true, $false, $false, $false, $false, $true, 1, $false, "72-45=27", 2) | Out-Null
$d.Content.Find.Execute("33+34=67", $true, $false, $false, $false, $false, $true, 1, $false, "42+8=50", 2) | Out-Null
$d.Content.Find.Execute("71-58=13", $true, $false, $false, $false, $false, $true, 1, $false, "5+24=29", 2) | Out-Null
$d.Content.Find.Execute("80+19=99", $true, $false, $false, $false, $false, $true, 1, $false, "57-13=44", 2) | Out-Null
$d.Content.Find.Execute("62-6=56", $true, $false, $false, $false, $false, $true, 1, $false, "63+25=88", 2) | Out-Null
$d.Content.Find.Execute("73-0=73", $true, $false, $false, $false, $false, $true, 1, $false, "35+11=46", 2) | Out-Null
$d.Content.Find.Execute("56-54=2", $true, $false, $false, $false, $false, $true, 1, $false, "63+21=84", 2) | Out-Null
$d.Content.Find.Execute("5+23=28", $true, $false, $false, $false, $false, $true, 1, $false, "42+56=98", 2) | Out-Null
$d.Content.Find.Execute("0+82=82", $true, $false, $false, $false, $false, $true, 1, $false, "36+35=71", 2) | Out-Null
$d.Content.Find.Execute("69+6=75", $true, $false, $false, $false, $false, $true, 1, $false, "1+45=46", 2) | Out-Null
$d.Content.Find.Execute("17-3=14", $true, $false, $false, $false, $false, $true, 1, $false, "4+11=15", 2) | Out-Null
$d.Content.Find.Execute("57-34=23", $true, $false, $false, $false, $false, $true, 1, $false, "83-57=26", 2) | Out-Null
$d.Content.Find.Execute("25+21=46", $true, $false, $false, $false, $false, $true, 1, $false, "73-15=58", 2) | Out-Null
$d.Content.Find.Execute("83-16=67", $true, $false, $false, $false, $false, $true, 1, $false, "76-37=39", 2) | Out-Null
$d.Content.Find.Execute("6+38=44", $true, $false, $false, $false, $false, $true, 1, $false, "6+14=20", 2) | Out-Null
$d.Content.Find.Execute("67-4=63", $true, $false, $false, $false, $false, $true, 1, $false, "23+5=28", 2) | Out-Null
$d.Content.Find.Execute("53+0=53", $true, $false, $false, $false, $false, $true, 1, $false, "82-5=77", 2) | Out-Null
$d.Content.Find.Execute("85-18=67", $true, $false, $false, $false, $false, $true, 1, $false, "0+39=39", 2) | Out-Null
$d.Content.Find.Execute("91-53=38", $true, $false, $false, $false, $false, $true, 1, $false, "18+57=75", 2) | Out-Null
$d.Content.Find.Execute("62-41=21", $true, $false, $false, $false, $false, $true, 1, $false, "66+22=88", 2) | Out-Null
$d.Content.Find.Execute("85+9=94", $true, $false, $false, $false, $false, $true, 1, $false, "99-44=55", 2) | Out-Null
$d.Content.Find.Execute("68+19=87", $true, $false, $false, $false, $false, $true, 1, $false, "11+36=47", 2) | Out-Null
$d.Content.Find.Execute("88-13=75", $true, $false, $false, $false, $false, $true, 1, $false, "60-48=12", 2) | Out-Null
$d.Content.Find.Execute("59-31=28", $true, $false, $false, $false, $false, $true, 1, $false, "19+25=44", 2) | Out-Null
$d.Content.Find.Execute("67-43=24", $true, $false, $false, $false, $false, $true, 1, $false, "55-3=52", 2) | Out-Null
$d.Content.Find.Execute("84-82=2", $true, $false, $false, $false, $false, $true, 1, $false, "49-2=47", 2) | Out-Null
$d.Content.Find.Execute("33+59=92", $true, $false, $false, $false, $false, $true, 1, $false, "68-27=41", 2) | Out-Null
$d.Content.Find.Execute("55-8=47", $true, $false, $false, $false, $false, $true, 1, $false, "30+51=81", 2) | Out-Null
$d.Content.Find.Execute("17-5=12", $true, $false, $false, $false, $false, $true, 1, $false, "8+22=30", 2) | Out-Null
$d.Content.Find.Execute("73+21=94", $true, $false, $false, $false, $false, $true, 1, $false, "3+85=88", 2) | Out-Null
$d.Content.Find.Execute("64+11=75", $true, $false, $false, $false, $false, $true, 1, $false, "93-46=47", 2) | Out-Null
$d.Content.Find.Execute("48+34=82", $true, $false, $false, $false, $false, $true, 1, $false, "87-54=33", 2) | Out-Null
$d.Content.Find.Execute("43+47=90", $true, $false, $false, $false, $false, $true, 1, $false, "90-60=30", 2) | Out-Null
$d.Content.Find.Execute("9+68=77", $true, $false, $false, $false, $false, $true, 1, $false, "60+17=77", 2) | Out-Null
$d.Content.Find.Execute("78+21=99", $true, $false, $false, $false, $false, $true, 1, $false, "58-23=35", 2) | Out-Null
$d.Content.Find.Execute("80-66=14", $true, $false, $false, $false, $false, $true, 1, $false, "29+42=71", 2) | Out-Null
$d.Content.Find.Execute("75-71=4", $true, $false, $false, $false, $false, $true, 1, $false, "52+4=56", 2) | Out-Null
$d.Content.Find.Execute("30+27=57", $true, $false, $false, $false, $false, $true, 1, $false, "16-16=0", 2) | Out-Null
$d.Content.Find.Execute("23+9=32", $true, $false, $false, $false, $false, $true, 1, $false, "6+61=67", 2) | Out-Null
$d.Content.Find.Execute("27-2=25", $true, $false, $false, $false, $false, $true, 1, $false, "8+78=86", 2) | Out-Null
$d.Content.Find.Execute("74-38=36", $true, $false, $false, $false, $false, $true, 1, $false, "37+39=76", 2) | Out-Null
$d.Content.Find.Execute("26+23=49", $true, $false, $false, $false, $false, $true, 1, $false, "3+16=19", 2) | Out-Null
$d.Content.Find.Execute("40+39=79", $true, $false, $false, $false, $false, $true, 1, $false, "94-47=47", 2) | Out-Null
$d.Content.Find.Execute("78-53=25", $true, $false, $false, $false, $false, $true, 1, $false, "30-19=11", 2) | Out-Null
$d.Content.Find.Execute("12+60=72", $true, $false, $false, $false, $false, $true, 1, $false, "48-31=17", 2) | Out-Null
$d.Content.Find.Execute("66-25=41", $true, $false, $false, $false, $false, $true, 1, $false, "12+18=30", 2) | Out-Null
$d.Content.Find.Execute("52-49=3", $true, $false, $false, $false, $false, $true, 1, $false, "44-18=26", 2) | Out-Null
$d.Content.Find.Execute("67-8=59", $true, $false, $false, $false, $false, $true, 1, $false, "63-20=43", 2) | Out-Null
$d.Content.Find.Execute("99-21=78", $true, $false, $false, $false, $false, $true, 1, $false, "42-6=36", 2) | Out-Null
$d.Content.Find.Execute("6+29=35", $true, $false, $false, $false, $false, $true, 1, $false, "17-8=9", 2) | Out-Null
$d.Content.Find.Execute("73-61=12", $true, $false, $false, $false, $false, $true, 1, $false, "31+41=72", 2) | Out-Null
$d.Content.Find.Execute("91-41=50", $true, $false, $false, $false, $false, $true, 1, $false, "86-81=5", 2) | Out-Null
$d.Content.Find.Execute("3+75=78", $true, $false, $false, $false, $false, $true, 1, $false, "35-17=18", 2) | Out-Null
$d.Content.Find.Execute("94-93=1", $true, $false, $false, $false, $false, $true, 1, $false, "3+59=62", 2) | Out-Null
$d.Content.Find.Execute("60-20=40", $true, $false, $false, $false, $false, $true, 1, $false, "14+2=16", 2) | Out-Null
$d.Content.Find.Execute("10+48=58", $true, $false, $false, $false, $false, $true, 1, $false, "90-56=34", 2) | Out-Null
$d.Content.Find.Execute("26+69=95", $true, $false, $false, $false, $false, $true, 1, $false, "9+9=18", 2) | Out-Null
$d.Content.Find.Execute("78+0=78", $true, $false, $false, $false, $false, $true, 1, $false, "87-80=7", 2) | Out-Null
$d.Content.Find.Execute("77+14=91", $true, $false, $false, $false, $false, $true, 1, $false, "65-33=32", 2) | Out-Null
$d.Content.Find.Execute("20+0=20", $true, $false, $false, $false, $false, $true, 1, $false, "99-65=34", 2) | Out-Null
$d.Content.Find.Execute("79-9=70", $true, $false, $false, $false, $false, $true, 1, $false, "97-76=21", 2) | Out-Null
$d.Content.Find.Execute("95-91=4", $true, $false, $false, $false, $false, $true, 1, $false, "31+5=36", 2) | Out-Null
$d.Content.Find.Execute("48-17=31", $true, $false, $false, $false, $false, $true, 1, $false, "89-32=57", 2) | Out-Null
$d.Content.Find.Execute("88-14=74", $true, $false, $false, $false, $false, $true, 1, $false, "42+11=53", 2) | Out-Null
$d.Content.Find.Execute("4+54=58", $true, $false, $false, $false, $false, $true, 1, $false, "81+6=87", 2) | Out-Null
$d.Content.Find.Execute("83-63=20", $true, $false, $false, $false, $false, $true, 1, $false, "61-29=32", 2) | Out-Null
$d.Content.Find.Execute("76-2=74", $true, $false, $false, $false, $false, $true, 1, $false, "87-23=64", 2) | Out-Null
$d.Content.Find.Execute("42+2=44", $true, $false, $false, $false, $false, $true, 1, $false, "38+26=64", 2) | Out-Null
$d.Content.Find.Execute("12+10=22", $true, $false, $false, $false, $false, $true, 1, $false, "97-83=14", 2) | Out-Null
$d.Content.Find.Execute("71-50=21", $true, $false, $false, $false, $false, $true, 1, $false, "53-31=22", 2) | Out-Null
$d.Content.Find.Execute("36+13=49", $true, $false, $false, $false, $false, $true, 1, $false, "13+8=21", 2) | Out-Null
$d.Content.Find.Execute("82-76=6", $true, $false, $false, $false, $false, $true, 1, $false, "6+72=78", 2) | Out-Null
$d.Content.Find.Execute("30+38=68", $true, $false, $false, $false, $false, $true, 1, $false, "79-3=76", 2) | Out-Null
$d.Content.Find.Execute("11+80=91", $true, $false, $false, $false, $false, $true, 1, $false, "89-20=69", 2) | Out-Null
$d.Content.Find.Execute("56+26=82", $true, $false, $false, $false, $false, $true, 1, $false, "24+4=28", 2) | Out-Null
$d.Content.Find.Execute("12-0=12", $true, $false, $false, $false, $false, $true, 1, $false, "83-79=4", 2) | Out-Null
$d.Content.Find.Execute("48+41=89", $true, $false, $false, $false, $false, $true, 1, $false, "41-15=26", 2) | Out-Null
